$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CEA")

$ws.Range("D8").Value = 15208300
$ws.Range("E8").Value = 14678300
$ws.Range("F8").Value = 13945900
$ws.Range("G8").Value = 13384400
$ws.Range("H8").Value = 13096400
$ws.Range("I8").Value = 12652400
$ws.Range("J8").Value = 12229400
$ws.Range("D9").Value = 10446400
$ws.Range("E9").Value = 9339400
$ws.Range("F9").Value = 8885900
$ws.Range("G9").Value = 9314700
$ws.Range("H9").Value = 9704300
$ws.Range("I9").Value = 8319000
$ws.Range("J9").Value = 7901300
$ws.Range("D10").Value = 4761900
$ws.Range("E10").Value = 5338900
$ws.Range("F10").Value = 5060000
$ws.Range("G10").Value = 4069700
$ws.Range("H10").Value = 3392200
$ws.Range("I10").Value = 4333400
$ws.Range("J10").Value = 4328100
$ws.Range("D14").Value = -206100
$ws.Range("E14").Value = 2100
$ws.Range("F14").Value = 26900
$ws.Range("J14").Value = 94700
$ws.Range("D15").Value = 2073100
$ws.Range("E15").Value = 1803800
$ws.Range("F15").Value = 1554000
$ws.Range("G15").Value = 1362800
$ws.Range("H15").Value = 1220900
$ws.Range("I15").Value = 1121500
$ws.Range("J15").Value = 1033800
$ws.Range("D17").Value = 13808700
$ws.Range("E17").Value = 12825300
$ws.Range("F17").Value = 12072300
$ws.Range("G17").Value = 12485300
$ws.Range("H17").Value = 12862500
$ws.Range("I17").Value = 12025000
$ws.Range("J17").Value = 11610200
$ws.Range("D18").Value = 1399700
$ws.Range("E18").Value = 1853000
$ws.Range("F18").Value = 1873700
$ws.Range("G18").Value = 899100
$ws.Range("H18").Value = 233900
$ws.Range("I18").Value = 627500
$ws.Range("J18").Value = 619200
$ws.Range("D20").Value = 350700
$ws.Range("E20").Value = -488300
$ws.Range("F20").Value = -707800
$ws.Range("G20").Value = 1800
$ws.Range("H20").Value = 324900
$ws.Range("I20").Value = 71500
$ws.Range("J20").Value = 316200
$ws.Range("D21").Value = 3820400
$ws.Range("E21").Value = 3220900
$ws.Range("F21").Value = 2771400
$ws.Range("G21").Value = 2257500
$ws.Range("H21").Value = 1774100
$ws.Range("I21").Value = 1815400
$ws.Range("J21").Value = 1963100
$ws.Range("D22").Value = 472500
$ws.Range("E22").Value = 400600
$ws.Range("F22").Value = 324900
$ws.Range("G22").Value = 290400
$ws.Range("H22").Value = 229800
$ws.Range("I22").Value = 251900
$ws.Range("J22").Value = 217100
$ws.Range("D23").Value = 1277800
$ws.Range("E23").Value = 964200
$ws.Range("F23").Value = 841000
$ws.Range("G23").Value = 610400
$ws.Range("H23").Value = 329000
$ws.Range("I23").Value = 447100
$ws.Range("J23").Value = 718400
$ws.Range("D24").Value = 267100
$ws.Range("E24").Value = 228800
$ws.Range("F24").Value = 92600
$ws.Range("G24").Value = 85000
$ws.Range("H24").Value = 18400
$ws.Range("I24").Value = 30400
$ws.Range("J24").Value = 39200
$ws.Range("D26").Value = 1010700
$ws.Range("E26").Value = 735400
$ws.Range("F26").Value = 748400
$ws.Range("G26").Value = 525400
$ws.Range("H26").Value = 310600
$ws.Range("I26").Value = 416700
$ws.Range("J26").Value = 679200
$ws.Range("D27").Value = 941200
$ws.Range("E27").Value = 667500
$ws.Range("F27").Value = 673300
$ws.Range("G27").Value = 506100
$ws.Range("H27").Value = 352100
$ws.Range("I27").Value = 438400
$ws.Range("J27").Value = 679100
$ws.Range("D32").Value = -350700
$ws.Range("E32").Value = 488300
$ws.Range("F32").Value = 707800
$ws.Range("G32").Value = -1800
$ws.Range("H32").Value = -324900
$ws.Range("I32").Value = -71500
$ws.Range("J32").Value = -316200
$ws.Range("D33").Value = 941200
$ws.Range("E33").Value = 667500
$ws.Range("F33").Value = 673300
$ws.Range("G33").Value = 506100
$ws.Range("H33").Value = 352100
$ws.Range("I33").Value = 438400
$ws.Range("J33").Value = 679100
$ws.Range("D35").Value = 941200
$ws.Range("E35").Value = 667500
$ws.Range("F35").Value = 673300
$ws.Range("G35").Value = 506100
$ws.Range("H35").Value = 352100
$ws.Range("I35").Value = 438400
$ws.Range("J35").Value = 679100
$ws.Range("D41").Value = 683400
$ws.Range("E41").Value = 251600
$ws.Range("F41").Value = 1347600
$ws.Range("G41").Value = 201100
$ws.Range("H41").Value = 296100
$ws.Range("I41").Value = 372800
$ws.Range("J41").Value = 573000
$ws.Range("D43").Value = 1549700
$ws.Range("E43").Value = 1307300
$ws.Range("F43").Value = 1505300
$ws.Range("G43").Value = 573200
$ws.Range("H43").Value = 523100
$ws.Range("I43").Value = 439600
$ws.Range("J43").Value = 371600
$ws.Range("D44").Value = 324300
$ws.Range("E44").Value = 333600
$ws.Range("F44").Value = 305100
$ws.Range("G44").Value = 335300
$ws.Range("H44").Value = 342100
$ws.Range("I44").Value = 309900
$ws.Range("J44").Value = 230900
$ws.Range("D45").Value = 157500
$ws.Range("E45").Value = 465400
$ws.Range("F45").Value = 267000
$ws.Range("G45").Value = 1597900
$ws.Range("H45").Value = 710100
$ws.Range("I45").Value = 758800
$ws.Range("J45").Value = 858900
$ws.Range("D46").Value = 2714900
$ws.Range("E46").Value = 2357900
$ws.Range("F46").Value = 3425000
$ws.Range("G46").Value = 2707400
$ws.Range("H46").Value = 1871400
$ws.Range("I46").Value = 1881100
$ws.Range("J46").Value = 2034400
$ws.Range("D47").Value = 456500
$ws.Range("E47").Value = 413800
$ws.Range("F47").Value = 517500
$ws.Range("G47").Value = 300400
$ws.Range("H47").Value = 283200
$ws.Range("I47").Value = 220600
$ws.Range("J47").Value = 222800
$ws.Range("D48").Value = 28513600
$ws.Range("E48").Value = 26289800
$ws.Range("F48").Value = 23015600
$ws.Range("G48").Value = 19248600
$ws.Range("H48").Value = 16188400
$ws.Range("I48").Value = 14011900
$ws.Range("J48").Value = 12574200
$ws.Range("D49").Value = 1721000
$ws.Range("E49").Value = 1725100
$ws.Range("F49").Value = 1710000
$ws.Range("G49").Value = 1706700
$ws.Range("H49").Value = 1705200
$ws.Range("I49").Value = 1699200
$ws.Range("J49").Value = 1685000
$ws.Range("D52").Value = 687900
$ws.Range("E52").Value = 724400
$ws.Range("F52").Value = 715900
$ws.Range("G52").Value = 647500
$ws.Range("H52").Value = 739200
$ws.Range("I52").Value = 563200
$ws.Range("J52").Value = 512000
$ws.Range("D54").Value = 34093800
$ws.Range("E54").Value = 31511000
$ws.Range("F54").Value = 29384000
$ws.Range("G54").Value = 24610700
$ws.Range("H54").Value = 20787500
$ws.Range("I54").Value = 18375900
$ws.Range("J54").Value = 17028400
$ws.Range("D57").Value = 472500
$ws.Range("E57").Value = 501000
$ws.Range("F57").Value = 550900
$ws.Range("G57").Value = 309100
$ws.Range("H57").Value = 3207000
$ws.Range("I57").Value = 485300
$ws.Range("J57").Value = 413900
$ws.Range("D58").Value = 7172800
$ws.Range("E58").Value = 5237200
$ws.Range("F58").Value = 6578000
$ws.Range("G58").Value = 4937900
$ws.Range("H58").Value = 3898100
$ws.Range("I58").Value = 3855600
$ws.Range("J58").Value = 3117500
$ws.Range("D59").Value = 4276100
$ws.Range("E59").Value = 4365800
$ws.Range("F59").Value = 3910900
$ws.Range("G59").Value = 3825300
$ws.Range("H59").Value = 772800
$ws.Range("I59").Value = 2875200
$ws.Range("J59").Value = 2907600
$ws.Range("D60").Value = 11921500
$ws.Range("E60").Value = 10104000
$ws.Range("F60").Value = 11039800
$ws.Range("G60").Value = 9072300
$ws.Range("H60").Value = 7877900
$ws.Range("I60").Value = 7216100
$ws.Range("J60").Value = 6439000
$ws.Range("D61").Value = 12219800
$ws.Range("E61").Value = 12241500
$ws.Range("F61").Value = 11099300
$ws.Range("G61").Value = 9589100
$ws.Range("H61").Value = 7045000
$ws.Range("I61").Value = 6285000
$ws.Range("J61").Value = 6144900
$ws.Range("D62").Value = 1229300
$ws.Range("E62").Value = 1393900
$ws.Range("F62").Value = 1318800
$ws.Range("G62").Value = 1234200
$ws.Range("H62").Value = 1622700
$ws.Range("I62").Value = 1230800
$ws.Range("J62").Value = 1208400
$ws.Range("D66").Value = 25877800
$ws.Range("E66").Value = 24172100
$ws.Range("F66").Value = 23831800
$ws.Range("G66").Value = 20162200
$ws.Range("H66").Value = 16794900
$ws.Range("I66").Value = 14973500
$ws.Range("J66").Value = 14041500
$ws.Range("D72").Value = 6068900
$ws.Range("E72").Value = 5191800
$ws.Range("F72").Value = 3602100
$ws.Range("G72").Value = -419900
$ws.Range("H72").Value = -877300
$ws.Range("I72").Value = -944700
$ws.Range("J72").Value = -1382000
$ws.Range("D76").Value = 8216000
$ws.Range("E76").Value = 7338900
$ws.Range("F76").Value = 5552200
$ws.Range("G76").Value = 4448400
$ws.Range("H76").Value = 3992600
$ws.Range("I76").Value = 3402400
$ws.Range("J76").Value = 2986900
$ws.Range("D81").Value = 941200
$ws.Range("E81").Value = 667500
$ws.Range("F81").Value = 673300
$ws.Range("G81").Value = 506100
$ws.Range("H81").Value = 352100
$ws.Range("I81").Value = 438400
$ws.Range("J81").Value = 679100
$ws.Range("D83").Value = 2066500
$ws.Range("E83").Value = 1852900
$ws.Range("F83").Value = 1602700
$ws.Range("G83").Value = 1354200
$ws.Range("H83").Value = 1213100
$ws.Range("I83").Value = 1114400
$ws.Range("J83").Value = 1025800
$ws.Range("D89").Value = 2904700
$ws.Range("E89").Value = 3694400
$ws.Range("F89").Value = 3610100
$ws.Range("G89").Value = 1824800
$ws.Range("H89").Value = 1603800
$ws.Range("I89").Value = 1872500
$ws.Range("J89").Value = 2021800
$ws.Range("D91").Value = -1157000
$ws.Range("E91").Value = -3195700
$ws.Range("F91").Value = -1277700
$ws.Range("G91").Value = -864900
$ws.Range("H91").Value = -270400
$ws.Range("I91").Value = -912400
$ws.Range("J91").Value = -796700
$ws.Range("D94").Value = -3162900
$ws.Range("E94").Value = -5517900
$ws.Range("F94").Value = -4125800
$ws.Range("G94").Value = -3566700
$ws.Range("H94").Value = -2527100
$ws.Range("I94").Value = -1749600
$ws.Range("J94").Value = -2217100
$ws.Range("D96").Value = -105200
$ws.Range("E96").Value = -109500
$ws.Range("D100").Value = 698700
$ws.Range("E100").Value = 687700
$ws.Range("F100").Value = 1644800
$ws.Range("G100").Value = 1649100
$ws.Range("H100").Value = 850300
$ws.Range("I100").Value = -322700
$ws.Range("J100").Value = 316900
$ws.Range("D101").Value = -7000
$ws.Range("E101").Value = 39800
$ws.Range("F101").Value = 17400
$ws.Range("H101").Value = -3800
$ws.Range("J101").Value = -5500
$ws.Range("D102").Value = 433500
$ws.Range("E102").Value = -1096000
$ws.Range("F102").Value = 1146500
$ws.Range("G102").Value = -95000
$ws.Range("H102").Value = -76700
$ws.Range("I102").Value = -200200
$ws.Range("J102").Value = 116200
